# Add the new "Wind Sensor Measurements" section (rows 30-35) to the
# InflowWind worksheet, matching the added lidar-related output channels
# (WindMeas1..WindMeas5).

$wb = $excel.ActiveWorkbook
$wsInstructions = $wb.Worksheets.Item("Instructions")
$ws = $wb.Worksheets.Item("InflowWind")

# First copy the formatting down (so the new cells inherit the correct
# styles), then fill in the values. Copy source ranges are split so that
# columns with no content in the source row (C on row 29, B on row 2)
# don't materialize stray empty cells in the destination.

# Rows 31-35 formats, from row 29 (the last existing data row).
$ws.Range("B29").Copy() | Out-Null
$ws.Range("B31:B35").PasteSpecial(-4122) | Out-Null
$ws.Range("D29:F29").Copy() | Out-Null
$ws.Range("D31:F35").PasteSpecial(-4122) | Out-Null

# Row 30 format, from row 2 (the existing category header row).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:E2").Copy() | Out-Null
$ws.Range("C30:E30").PasteSpecial(-4122) | Out-Null
# The old placeholder row 30 had a leftover styled F30 cell; remove it
# since the new category-header row has no F column content/style.
$ws.Range("F30").Clear() | Out-Null

$names = @("WindMeas1", "WindMeas2", "WindMeas3", "WindMeas4", "WindMeas5")
$descriptions = @(
    "Wind measurement at sensor 1",
    "Wind measurement at sensor 2",
    "Wind measurement at sensor 3",
    "Wind measurement at sensor 4",
    "Wind measurement at sensor 5"
)
$criteria = @(
    "p%lidar%SensorType == SensorType_None",
    "p%lidar%NumPulseGate < 2",
    "p%lidar%NumPulseGate < 3",
    "p%lidar%NumPulseGate < 4",
    "p%lidar%NumPulseGate < 5"
)

# Column B (names) for all five rows first ...
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item(31 + $i, 2).Value = $names[$i]
}
# ... then column D (descriptions) for all five rows ...
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item(31 + $i, 4).Value = $descriptions[$i]
}
# ... then column E (convention) - same text for every row ...
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item(31 + $i, 5).Value = "Defined by sensor"
}
# ... then column F (units) - same text for every row ...
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item(31 + $i, 6).Value = "(m/s)"
}
# ... then the new category header in A30 ...
$ws.Range("A30").Value = "Wind Sensor Measurements"
# ... and finally column G (invalid channel criteria) for all five rows.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item(31 + $i, 7).Value = $criteria[$i]
}

# --- Update the selections to match the edit ---
$wsInstructions.Range("D8").Select() | Out-Null
$ws.Activate() | Out-Null
$ws.Range("B31").Select() | Out-Null

Write-Output "done"
